$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the placeholder "Your Name" / "Your Email" row with the real
# submitter's info, turning the email into a mailto: hyperlink just like
# the other two rows (John Doe / Jane Doe).
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:h.m.alrahmani@gmail.com", "", "", "h.m.alrahmani@gmail.com")
$ws.Range("A1").Value = "Hessa Hanim"
$ws.Range("B1").HorizontalAlignment = -4108

# Update the view state (zoom + selection) to match where the author left
# the workbook.
$excel.ActiveWindow.Zoom = 205
$null = $ws.Range("A5").Select()
